$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("D6").Value = "PC4 UART7 Rx"
$ws.Range("D7").Value = "PC5 UART7 Tx"
$ws.Range("D8").Value = "PC6"
$ws.Range("D9").Value = "PE5"

$ws.Activate()
$ws.Range("D34").Select()
